# Regenerate the handback/localization status report for archive:
#  - The "Status" column value flips from "Ready for handoff" to
#    "In Translation" everywhere it appears (Overview!E2:F2, zh-cn!C2,
#    de-de!C2 all share that string).
#  - With the shorter text in place the "Status" column is narrower, so
#    its column width shrinks to match on all three sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Update the status text (Overview tracks both locales in columns E/F).
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$zhcn.Range("C2").Value = "In Translation"
$dede.Range("C2").Value = "In Translation"

# Narrow the Status columns to fit the new, shorter text.
$overview.Columns.Item(5).ColumnWidth = 12.43
$overview.Columns.Item(6).ColumnWidth = 12.43
$zhcn.Columns.Item(3).ColumnWidth = 12.43
$dede.Columns.Item(3).ColumnWidth = 12.43
